$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A11").Value = "El usuario reporta dolor de cabeza y mareos desde hace dos días. Se recomienda consultar con un neurólogo debido a la posible relación con problemas neurológicos. `n"

$ws.Range("A19").Value = "El usuario mencionó que tiene dolor de cabeza y mareos desde hace dos días. `n"

$ws.Range("A26").Value = "Al llegar a la consulta, sería bueno que le preguntes al paciente qué otras cosas del cuerpo le molestan o si siente algún otro síntoma. Esto te dará una visión más amplia del problema y te ayudará a determinar si el dolor de cabeza y el mareo son parte de algo más grande. `n"
